$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 198 (pushes the existing rows 198-225 down to 199-226)
$ws.Rows.Item(198).Insert()

# Populate the newly inserted row 198 with the new weekly price record
$ws.Cells.Item(198, 1).Value = 11
$ws.Cells.Item(198, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(198, 3).Value = "Bíobío"
$ws.Cells.Item(198, 4).Value = 44504
$ws.Cells.Item(198, 5).Value = 8
$ws.Cells.Item(198, 6).Value = 100112006
$ws.Cells.Item(198, 7).Value = "Repollo"
$ws.Cells.Item(198, 8).Value = "Copenhague"
$ws.Cells.Item(198, 9).Value = "Primera"
$ws.Cells.Item(198, 10).Value = 1800
$ws.Cells.Item(198, 11).Value = 700
$ws.Cells.Item(198, 12).Value = 850
$ws.Cells.Item(198, 13).Value = 783
$ws.Cells.Item(198, 14).Value = "$/unidad"
$ws.Cells.Item(198, 15).Value = "Región del Maule"
$ws.Cells.Item(198, 16).Value = 783
$ws.Cells.Item(198, 17).Value = 1
$ws.Cells.Item(198, 18).Value = "Hortaliza"
